# Applies the "unitlibrary" edit: adds three new Unit Process rows
# (simple_DRI, simple_syngas, simple_casting) to the "Unit Processes" sheet,
# as part of cleaning up function calls / ensuring aggregate_flows is
# called properly, and setting up a testing file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data to append below the existing table (rows 44-46).
# Columns: A=ID  B=meta-category  C=display name  D=product  E=productType
#          F=varFile  G=varSheet  H=calcFile  I=calcSheet
$newRows = @(
    @{ Row = 44; A = "simple_DRI";     B = "steel";  C = "DRI Furnace";        D = "DRI";             E = "outflow"; F = "data/steel/steel_simplified_var.xlsx"; G = "DRI";    H = "data/steel/steel_simplified_calcs.xlsx"; I = "DRI" },
    @{ Row = 45; A = "simple_syngas";  B = "energy"; C = "Syngas Production";  D = "syngas";          E = "outflow"; F = "data/steel/steel_simplified_var.xlsx"; G = "syngas"; H = "data/steel/steel_simplified_calcs.xlsx"; I = "syngas" },
    @{ Row = 46; A = "simple_casting"; B = "steel";  C = "Finishing";          D = "hot rolled coil"; E = "outflow"; F = "data/steel/steel_simplified_var.xlsx"; G = "Finish"; H = "data/steel/steel_simplified_calcs.xlsx"; I = "Finish" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Columns A, C, D, E, F, G, H, I use a Text ("@") number format
    # (matching the existing rows' style); column B keeps the default style.
    $textCols = @("A", "C", "D", "E", "F", "G", "H", "I")
    foreach ($col in $textCols) {
        $ws.Range("$col$rowNum").NumberFormat = "@"
    }

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H
    $ws.Range("I$rowNum").Value = $r.I
}

# Replicate the author's view state change: freeze the header row and
# first column, then leave the view scrolled/selected near the bottom of
# the (now larger) table.
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("I49").Select() | Out-Null
